$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 798957.5600000001
$ws.Range("J17").Value = 798957.5600000001
$ws.Range("L17").Value = 2396872.68
$ws.Range("N17").Value = -2397208.68
$ws.Range("H40").Value = 2910.6428
$ws.Range("I40").Value = 5074.75
$ws.Range("J40").Value = 2045
$ws.Range("K40").Value = 5074.75
$ws.Range("L40").Value = 2045
$ws.Range("M40").Value = -4899.75
$ws.Range("N40").Value = -2395
$ws.Range("H43").Value = 18937.285
$ws.Range("I43").Value = 6245
$ws.Range("J43").Value = 35860.332
$ws.Range("K43").Value = 6245
$ws.Range("L43").Value = 35860.332
$ws.Range("M43").Value = -6176
$ws.Range("N43").Value = -35998.332
$ws.Range("H58").Value = 233.33333
$ws.Range("I58").Value = 100
$ws.Range("J58").Value = 500
$ws.Range("K58").Value = 300
$ws.Range("L58").Value = 1500
$ws.Range("M58").Value = -150
$ws.Range("N58").Value = -1800
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H112").Value = 1383.6735
$ws.Range("J112").Value = 1234.4445
$ws.Range("L112").Value = 3703.3335
$ws.Range("N112").Value = -5919.333500000001
$ws.Range("H115").Value = 453.8889
$ws.Range("I115").Value = 453.8889
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1361.6667
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 205.3333
$ws.Range("N115").ClearContents()
$ws.Range("H118").Value = 2143.375
$ws.Range("I118").Value = 1266.9231
$ws.Range("J118").Value = 2743.0527
$ws.Range("K118").Value = 3800.7693
$ws.Range("L118").Value = 8229.158100000001
$ws.Range("M118").Value = -2143.7693
$ws.Range("N118").Value = -11543.1581
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H133").Value = 47589.875
$ws.Range("J133").Value = 47589.875
$ws.Range("L133").Value = 47589.875
$ws.Range("N133").Value = -57709.875
$ws.Range("H136").Value = 50000
$ws.Range("J136").Value = 50000
$ws.Range("L136").Value = 50000
$ws.Range("N136").Value = -60200
$ws.Range("H138").Value = 5815924
$ws.Range("I138").Value = 1649.1515
$ws.Range("J138").Value = 25003032
$ws.Range("K138").Value = 4947.4545
$ws.Range("L138").Value = 75009096
$ws.Range("M138").Value = 192.5455000000002
$ws.Range("N138").Value = -75019376
$ws.Range("H140").Value = 39913.332
$ws.Range("I140").Value = 10000
$ws.Range("J140").Value = 54870
$ws.Range("K140").Value = 10000
$ws.Range("L140").Value = 54870
$ws.Range("M140").Value = -4820
$ws.Range("N140").Value = -65230

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1445.1666
$ws.Range("I2").Value = 1066.6666
$ws.Range("J2").Value = 1634.4166
$ws.Range("K2").Value = 1066.6666
$ws.Range("L2").Value = 1634.4166
$ws.Range("M2").Value = -953.6666
$ws.Range("N2").Value = -1860.4166
$ws.Range("H32").Value = 8717.575999999999
$ws.Range("I32").Value = 9019.121999999999
$ws.Range("J32").Value = 7240
$ws.Range("K32").Value = 9019.121999999999
$ws.Range("L32").Value = 7240
$ws.Range("M32").Value = -8732.121999999999
$ws.Range("N32").Value = -7814
$ws.Range("H45").Value = 2725.45
$ws.Range("I45").Value = 2800.6428
$ws.Range("J45").Value = 2550
$ws.Range("K45").Value = 2800.6428
$ws.Range("L45").Value = 2550
$ws.Range("M45").Value = -2423.6428
$ws.Range("N45").Value = -3304
$ws.Range("H61").Value = 9617488
$ws.Range("I61").Value = 13159698
$ws.Range("J61").Value = 2917.2856
$ws.Range("K61").Value = 13159698
$ws.Range("L61").Value = 2917.2856
$ws.Range("M61").Value = -13159486
$ws.Range("N61").Value = -3341.2856
$ws.Range("H74").Value = 8198448
$ws.Range("I74").Value = 11364462
$ws.Range("J74").Value = 4057.7646
$ws.Range("K74").Value = 11364462
$ws.Range("L74").Value = 4057.7646
$ws.Range("M74").Value = -11363588
$ws.Range("N74").Value = -5805.7646
$ws.Range("H77").Value = 8198448
$ws.Range("I77").Value = 11364462
$ws.Range("J77").Value = 4057.7646
$ws.Range("K77").Value = 56822310
$ws.Range("L77").Value = 20288.823
$ws.Range("M77").Value = -56817942
$ws.Range("N77").Value = -29024.823
$ws.Range("H116").Value = 1445.1666
$ws.Range("I116").Value = 1066.6666
$ws.Range("J116").Value = 1634.4166
$ws.Range("K116").Value = 1066.6666
$ws.Range("L116").Value = 1634.4166
$ws.Range("M116").Value = 1227.3334
$ws.Range("N116").Value = -6222.4166
$ws.Range("H132").Value = 5815782
$ws.Range("I132").Value = 7813759
$ws.Range("K132").Value = 23441277
$ws.Range("M132").Value = -23438747
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 9617488
$ws.Range("I136").Value = 13159698
$ws.Range("J136").Value = 2917.2856
$ws.Range("K136").Value = 39479094
$ws.Range("L136").Value = 8751.856800000001
$ws.Range("M136").Value = -39476544
$ws.Range("N136").Value = -13851.8568

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1445.1666
$ws.Range("I3").Value = 1066.6666
$ws.Range("J3").Value = 1634.4166
$ws.Range("K3").Value = 1066.6666
$ws.Range("L3").Value = 1634.4166
$ws.Range("M3").Value = -952.6666
$ws.Range("N3").Value = -1862.4166
$ws.Range("H20").Value = 935.24
$ws.Range("I20").Value = 944.6667
$ws.Range("J20").Value = 911
$ws.Range("K20").Value = 944.6667
$ws.Range("L20").Value = 911
$ws.Range("M20").Value = -697.6667
$ws.Range("N20").Value = -1405
$ws.Range("H107").Value = 2907.3125
$ws.Range("I107").Value = 2893.8572
$ws.Range("K107").Value = 2893.8572
$ws.Range("M107").Value = -973.8571999999999
$ws.Range("N107").ClearContents()
$ws.Range("H139").Value = 22198.75
$ws.Range("J139").Value = 22198.75
$ws.Range("L139").Value = 22198.75
$ws.Range("N139").Value = -32478.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2016.4286
$ws.Range("I122").Value = 2120.9092
$ws.Range("J122").Value = 1633.3334
$ws.Range("K122").Value = 6362.7276
$ws.Range("L122").Value = 4900.0002
$ws.Range("M122").Value = -3912.7276
$ws.Range("N122").Value = -9800.0002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 817.1429000000001
$ws.Range("I46").Value = 165.71428
$ws.Range("K46").Value = 497.14284
$ws.Range("M46").Value = -406.14284
$ws.Range("N46").ClearContents()
$ws.Range("H118").Value = 3063.2
$ws.Range("J118").Value = 3447.75
$ws.Range("L118").Value = 10343.25
$ws.Range("N118").Value = -12829.25
$ws.Range("H122").Value = 1178.1875
$ws.Range("I122").Value = 1706.8572
$ws.Range("J122").Value = 767
$ws.Range("K122").Value = 15361.7148
$ws.Range("L122").Value = 6903
$ws.Range("M122").Value = -12911.7148
$ws.Range("N122").Value = -11803
$ws.Range("H132").Value = 2135.0688
$ws.Range("J132").Value = 3131.8125
$ws.Range("L132").Value = 28186.3125
$ws.Range("N132").Value = -33246.3125

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3824.182
$ws.Range("I132").Value = 2159.611
$ws.Range("J132").Value = 5821.6665
$ws.Range("K132").Value = 6478.833
$ws.Range("L132").Value = 17464.9995
$ws.Range("M132").Value = -3948.833
$ws.Range("N132").Value = -22524.9995

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 952.2143
$ws.Range("I22").Value = 812.8570999999999
$ws.Range("J22").Value = 1091.5714
$ws.Range("K22").Value = 812.8570999999999
$ws.Range("L22").Value = 1091.5714
$ws.Range("M22").Value = -517.8570999999999
$ws.Range("N22").Value = -1681.5714
$ws.Range("H27").Value = 952.2143
$ws.Range("I27").Value = 812.8570999999999
$ws.Range("J27").Value = 1091.5714
$ws.Range("K27").Value = 812.8570999999999
$ws.Range("L27").Value = 1091.5714
$ws.Range("M27").Value = -705.8570999999999
$ws.Range("N27").Value = -1305.5714
$ws.Range("H46").Value = 1001.4
$ws.Range("I46").Value = 932
$ws.Range("J46").Value = 1070.8
$ws.Range("K46").Value = 932
$ws.Range("L46").Value = 1070.8
$ws.Range("M46").Value = -744
$ws.Range("N46").Value = -1446.8
$ws.Range("H122").Value = 6337.826
$ws.Range("I122").Value = 6415.357
$ws.Range("K122").Value = 19246.071
$ws.Range("M122").Value = -16796.071
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 35719260
$ws.Range("I136").Value = 62501816
$ws.Range("J136").Value = 9181.111000000001
$ws.Range("K136").Value = 187505448
$ws.Range("L136").Value = 27543.333
$ws.Range("M136").Value = -187502898
$ws.Range("N136").Value = -32643.333

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 26404.908
$ws.Range("J64").Value = 30050.445
$ws.Range("L64").Value = 30050.445
$ws.Range("N64").Value = -30546.445
$ws.Range("H67").Value = 26404.908
$ws.Range("J67").Value = 30050.445
$ws.Range("L67").Value = 30050.445
$ws.Range("N67").Value = -31766.445
$ws.Range("H122").Value = 2671.6365
$ws.Range("I122").Value = 2526.4443
$ws.Range("K122").Value = 7579.3329
$ws.Range("M122").Value = -5129.3329
$ws.Range("N122").ClearContents()
